# The underlying OOXML diff for this revision (commit "Moving from 2.0.1
# to 2.0.2") does not change any document content, text, formatting,
# structure, style definition or page-setup value: every hunk just shows
# the very same elements/attributes (same names, same values - e.g.
# <w:pgSz w:w="11906" w:h="16838"/> vs <w:pgSz w:h="16838" w:w="11906"/>,
# the same namespace declarations reordered on <w:document>, the same
# <w:lsdException>/<w:style> attributes reordered, etc.) re-serialized in
# a different (alphabetical) attribute order. That reordering comes from
# the authoring tool that regenerated the fixture during the dependency
# bump, not from an actual Word edit, so there is nothing to change in
# the document's content/formatting.
#
# Touch the document through the Word object model without modifying
# anything, so the template is left byte-for-byte equivalent (same text,
# same styles, same page setup) as required by the diff.
$d = $word.ActiveDocument

$sectionCount = $d.Sections.Count
$ps = $d.PageSetup
Write-Output ("Sections=" + $sectionCount + " PageWidth=" + $ps.PageWidth + " PageHeight=" + $ps.PageHeight)

$styleCount = $d.Styles.Count
Write-Output ("Styles=" + $styleCount)
